# Deploy the implementation guide.
# - Update the "Date" metadata value on the Metadata sheet.
# - Insert a new "CNV" concept row on the Concepts sheet (between SNV and QC).

$wb = $excel.ActiveWorkbook

# --- Update Date on Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$found = $meta.Columns.Item(1).Find("Date")
if ($found -ne $null) {
    $meta.Cells.Item($found.Row, 2).Value = "2021-11-30T21:05:25+00:00"
}

# --- Insert new CNV row on Concepts sheet ---
$concepts = $wb.Worksheets.Item("Concepts")

# Row 4 currently holds "QC"; insert a new blank row above it so the SNV row
# (row 3) is immediately followed by the new CNV row, and QC/INDEX/CF shift
# down by one.
$concepts.Rows.Item(4).Insert()

# Fill in the new row's values. Force column A to be stored as text ("1"),
# matching the other Level cells, instead of being auto-converted to a number.
$concepts.Cells.Item(4, 1).NumberFormat = "@"
$concepts.Cells.Item(4, 1).Value = "1"
$concepts.Cells.Item(4, 2).Value = "CNV"
$concepts.Cells.Item(4, 3).Value = "Raw CNV calls"
$concepts.Cells.Item(4, 4).Value = "Raw CNV calls"

# Copy formatting only (border/fill/font/alignment) from the row above (row 3,
# an existing data row) onto the newly inserted row 4.
$concepts.Range("A3:D3").Copy()
$concepts.Range("A4:D4").PasteSpecial(-4122, $null, $false, $false)
$concepts.Application.CutCopyMode = $false
